# Updated remaining queries for C3DC
# Applies the JOIN-column rename (std.id/prt.id -> std.study_id/prt.participant_id, etc.)
# to every SQL query cell on the sheet, adjusts the selected/scrolled cell,
# and widens column C to fit the (now slightly longer) query text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Update-Query([string]$t) {
    $t = $t.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $t = $t.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $t = $t.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $t = $t.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $t = $t.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $t = $t.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')
    return $t
}

# Every query cell on the sheet (Studies, Summary, Participants, Diagnosis,
# Treatment, TreatmentResp, Survival) shares the same stale JOIN clauses.
$cells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $cells) {
    $orig = $ws.Range($addr).Text
    $updated = Update-Query $orig
    if ($updated -ne $orig) {
        $ws.Range($addr).Value = $updated
    }
}

# Column C (the StatQuery column) now needs to be a bit wider since it no
# longer relies on auto "best fit" sizing.
$ws.Columns.Item(3).ColumnWidth = 67.3

# Selection / scroll moved down to the Treatment row while editing.
$ws.Range("C7").Select()
